# Updated loading-percent results for the 380 kV case (row index = A column value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row -> { column letter = new value } for the columns that changed
# (B, C, D, E, F, H, L, N). Columns G, I, J, K, M, O stay at 0 and are left untouched.
$updates = @{
    2 = @{ "B"=24.57465587734637; "C"=9.256536877775636; "D"=7.852815787582229; "E"=9.589500852679581; "F"=41.57193039689167; "H"=7.344005520526261; "L"=10.48871128386846; "N"=20.10905685934477 }
    3 = @{ "B"=24.07473069002777; "C"=8.639776147237788; "D"=7.875895845758047; "E"=9.608231250050473; "F"=41.2008612623715; "H"=7.344005520526261; "L"=10.47213174224119; "N"=20.16594502128077 }
    4 = @{ "B"=23.77074027902026; "C"=8.237740905599523; "D"=7.891315933408685; "E"=9.620485003206959; "F"=40.9843093299569; "H"=7.344005520526261; "L"=10.46427795857785; "N"=20.20286001205965 }
    5 = @{ "B"=23.64779423517601; "C"=8.068018890434301; "D"=7.897912594986595; "E"=9.625668289068519; "F"=40.89897459810155; "H"=7.344005520526261; "L"=10.46166356297643; "N"=20.21840202707692 }
    6 = @{ "B"=23.62744093717568; "C"=8.039480391074317; "D"=7.899026826828183; "E"=9.626540443350621; "F"=40.88498273707329; "H"=7.344005520526261; "L"=10.46126486766305; "N"=20.22101288555718 }
    7 = @{ "B"=23.76907816234892; "C"=8.235475847662768; "D"=7.891403633108447; "E"=9.620554137878608; "F"=40.98314659264262; "H"=7.344005520526261; "L"=10.46424032567515; "N"=20.20306759773546 }
    8 = @{ "B"=24.40177773232991; "C"=9.048710301975317; "D"=7.860513881002285; "E"=9.595802992537925; "F"=41.44168969227907; "H"=7.344005520526261; "L"=10.48251270127108; "N"=20.12825957886528 }
    9 = @{ "B"=25.65793799425828; "C"=10.45915852073305; "D"=7.809901798322071; "E"=9.553225437205752; "F"=42.42689102791093; "H"=7.344005520526261; "L"=10.53673197986329; "N"=19.9973322749885 }
    10 = @{ "B"=26.57971266673569; "C"=11.38431239429845; "D"=7.778865465389941; "E"=9.525553289586037; "F"=43.19801815052898; "H"=7.344005520526261; "L"=10.58765641627453; "N"=19.91077829836105 }
    11 = @{ "B"=26.99670009313969; "C"=11.78133141734197; "D"=7.766097839835222; "E"=9.513743412342714; "F"=43.55791197821222; "H"=7.344005520526261; "L"=10.61319629520944; "N"=19.87350149933702 }
    12 = @{ "B"=27.15410227271438; "C"=11.92826517657317; "D"=7.761458749946631; "E"=9.509382875959673; "F"=43.69539951388224; "H"=7.344005520526261; "L"=10.62320523640169; "N"=19.85968796943787 }
    13 = @{ "B"=27.12022777162954; "C"=11.89677166229859; "D"=7.76244913207368; "E"=9.510317036637538; "F"=43.66573728793939; "H"=7.344005520526261; "L"=10.62103468546539; "N"=19.86264949774952 }
    14 = @{ "B"=27.00966064485445; "C"=11.79348795262505; "D"=7.765712247859424; "E"=9.513382433378686; "F"=43.56919965180992; "H"=7.344005520526261; "L"=10.61401299556405; "N"=19.87235898672665 }
    15 = @{ "B"=26.94186494042982; "C"=11.72978034865748; "D"=7.767736534611047; "E"=9.515274602818019; "F"=43.51022113002914; "H"=7.344005520526261; "L"=10.60975584410452; "N"=19.87834573411007 }
    16 = @{ "B"=26.55240017975752; "C"=11.35788822344564; "D"=7.779727169317383; "E"=9.526340727687382; "F"=43.17467246422938; "H"=7.344005520526261; "L"=10.58603477912578; "N"=19.91325669267001 }
    17 = @{ "B"=26.31275740521278; "C"=11.12365259573329; "D"=7.787430084208665; "E"=9.533328558406476; "F"=42.97108161761194; "H"=7.344005520526261; "L"=10.5720881400301; "N"=19.9352111134027 }
    18 = @{ "B"=26.17471107605079; "C"=10.98668155982642; "D"=7.791987666561871; "E"=9.537421051517951; "F"=42.85484497165262; "H"=7.344005520526261; "L"=10.56429012241622; "N"=19.94803604585935 }
    19 = @{ "B"=26.12794009139926; "C"=10.93991937229946; "D"=7.793552565599031; "E"=9.538819294543153; "F"=42.81564071565763; "H"=7.344005520526261; "L"=10.56168838180633; "N"=19.95241221468436 }
    20 = @{ "B"=26.33829077207369; "C"=11.14881957749092; "D"=7.786596934019413; "E"=9.532577110228017; "F"=42.99266560792142; "H"=7.344005520526261; "L"=10.57354965584807; "N"=19.93285359903686 }
    21 = @{ "B"=27.04215179268834; "C"=11.823917241341; "D"=7.764748467714385; "E"=9.512479025516761; "F"=43.59752328345677; "H"=7.344005520526261; "L"=10.61606630812328; "N"=19.86949886090756 }
    22 = @{ "B"=27.49916959382473; "C"=12.24527275301516; "D"=7.75161080971212; "E"=9.499994159851525; "F"=43.9997977071782; "H"=7.344005520526261; "L"=10.64581876654325; "N"=19.82985583129746 }
    23 = @{ "B"=27.25557731152014; "C"=12.02219825100439; "D"=7.758517679910192; "E"=9.506598158576983; "F"=43.78449436641162; "H"=7.344005520526261; "L"=10.62976087263972; "N"=19.85085246293843 }
    24 = @{ "B"=26.32674797297618; "C"=11.13744876870864; "D"=7.786973199294095; "E"=9.532916606018313; "F"=42.98290495356527; "H"=7.344005520526261; "L"=10.57288821872596; "N"=19.93391879890898 }
    25 = @{ "B"=25.31757069663357; "C"=10.09719606646479; "D"=7.822519489547913; "E"=9.564108244325958; "F"=42.15167167374197; "H"=7.344005520526261; "L"=10.52010627326199; "N"=20.03106089177414 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
